$d = $word.ActiveDocument

# The second table in the document holds the student (alumno/a) roster.
# It currently has a header row plus two data rows; the row for
# "Díez Viñas Malena" (the last row) must be removed entirely, along
# with its routes/dates, since that student entry was deleted.
$table = $d.Tables.Item(2)
$table.Rows.Item($table.Rows.Count).Delete()
